$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # params_ranges
$ws2 = $wb.Worksheets.Item(2)   # constraints_db
$ws3 = $wb.Worksheets.Item(3)   # requirements_db

# ---------------------------------------------------------------------------
# Sheet1 (params_ranges): insert two new columns (B "is an array", C "array
# size") before the existing type/min/max/step/fixed value columns, then
# populate the parameter rows (p4, p5, p6, p1Array, p7).
# ---------------------------------------------------------------------------
$ws1.Range("B1:C1").EntireColumn.Insert()
$ws1.Range("B1").Value = "is an array"
$ws1.Range("C1").Value = "array size"

$ws1.Range("A2").Value = "p4"
$ws1.Range("B2").Value = $false
$ws1.Range("C2").Value = 0
$ws1.Range("D2").Value = "FIXED"
$ws1.Range("E2").Value = 0
$ws1.Range("F2").Value = 0
$ws1.Range("G2").Value = 0
$ws1.Range("H2").Value = 0.8

$ws1.Range("A3").Value = "p5"
$ws1.Range("B3").Value = $false
$ws1.Range("C3").Value = 0
$ws1.Range("D3").Value = "BOOLEAN"
$ws1.Range("E3").Value = 0
$ws1.Range("F3").Value = 0
$ws1.Range("G3").Value = 0
$ws1.Range("H3").Value = 0

$ws1.Range("A4").Value = "p6"
$ws1.Range("B4").Value = $false
$ws1.Range("C4").Value = 0
$ws1.Range("D4").Value = "CONTINUOUS"
$ws1.Range("E4").Value = 0.1
$ws1.Range("F4").Value = 1
$ws1.Range("G4").Value = 0
$ws1.Range("H4").Value = 0

$ws1.Range("A5").Value = "p1Array"
$ws1.Range("B5").Value = $true
$ws1.Range("C5").Value = 3
$ws1.Range("D5").Value = "CONTINUOUS"
$ws1.Range("E5").Value = 1
$ws1.Range("F5").Value = 10
$ws1.Range("G5").Value = 0
$ws1.Range("H5").Value = 0
$ws1.Range("H5").HorizontalAlignment = -4152   # xlRight

$ws1.Range("A6").Value = "p7"
$ws1.Range("B6").Value = $false
$ws1.Range("C6").Value = 0
$ws1.Range("D6").Value = "DISCRETE"
$ws1.Range("E6").Value = 1
$ws1.Range("F6").Value = 3
$ws1.Range("G6").Value = 1
$ws1.Range("H6").Value = 0

# Data validations: the pre-existing CONTINUOUS/DISCRETE/FIXED/BOOLEAN list
# now lives on column D (it was auto-shifted by the column insert above);
# drop it and re-add so it ends up first, then add the new TRUE/FALSE
# validation for column B ("is an array").
$ws1.Range("D2:D1048576").Validation.Delete()
$ws1.Range("D2:D1048576").Validation.Add(3, 1, 1, """CONTINUOUS, DISCRETE, FIXED, BOOLEAN""")
$ws1.Range("B2:B1048576").Validation.Add(3, 1, 1, """TRUE, FALSE""")

$ws1.Columns.Item(1).ColumnWidth = 12.2
$ws1.Columns.Item(2).ColumnWidth = 12.2
$ws1.Columns.Item(3).ColumnWidth = 12.2
$ws1.Columns.Item(4).ColumnWidth = 14.2
$ws1.Columns.Item(5).ColumnWidth = 11.2
$ws1.Columns.Item(6).ColumnWidth = 11.8
$ws1.Columns.Item(8).ColumnWidth = 10.9

$ws1.Range("D11").Select()

# ---------------------------------------------------------------------------
# Sheet2 (constraints_db): insert a new column C "on array" between
# "enabled" and "expression", then populate the constraint rows.
# ---------------------------------------------------------------------------
$ws2.Range("C1").EntireColumn.Insert()
$ws2.Range("C1").Value = "on array"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = $true
$ws2.Range("C2").Value = $false
$ws2.Range("D2").Value = "p6+p7"
$ws2.Range("E2").Value = "<="
$ws2.Range("F2").Value = 2

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = $false
$ws2.Range("C3").Value = $false
$ws2.Range("D3").Value = "p2*p3"
$ws2.Range("E3").Value = ">="
$ws2.Range("F3").Value = 10

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = $true
$ws2.Range("C4").Value = $true
$ws2.Range("D4").Value = "2*p1Array[0]+p1Array[2]/2"
$ws2.Range("E4").Value = ">="
$ws2.Range("F4").Value = 8

$ws2.Range("B2:B1048576").Validation.Add(3, 1, 1, """TRUE, FALSE""")
$ws2.Range("C2:C1048576").Validation.Add(3, 1, 1, """TRUE, FALSE""")

$ws2.Columns.Item(1).ColumnWidth = 12.2
$ws2.Columns.Item(2).ColumnWidth = 13.2
$ws2.Columns.Item(3).ColumnWidth = 13.2
$ws2.Columns.Item(4).ColumnWidth = 22.7

$ws2.Range("F12").Select()

# ---------------------------------------------------------------------------
# Sheet3 (requirements_db): populate the first requirement row. The "type"
# cell holds ">=" typed with a leading apostrophe (quote-prefix) like the
# source workbook.
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = $false
$ws3.Range("C2").Value = "variable-variable1"
$ws3.Range("D2").Value = "'>="
$ws3.Range("E2").Value = 0

$ws3.Range("B2:B1048576").Validation.Add(3, 1, 1, """TRUE, FALSE""")

$ws3.Range("F10").Select()

# ---------------------------------------------------------------------------
# Workbook-level: requirements_db becomes the active/selected sheet.
# ---------------------------------------------------------------------------
$ws3.Activate()
